$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Outputs")
$ws.Activate()

# Decrement the "Run" index column (B2:B17) by 1, i.e. 1..16 -> 0..15
for ($row = 2; $row -le 17; $row++) {
    $cell = $ws.Cells.Item($row, 2)
    $current = $cell.Value2
    $cell.Value = $current - 1
}

# Update the active selection on the sheet to C2
$ws.Range("C2").Select()
